# Updates the cryptocurrency price/volume snapshot table on Sheet1
# (columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)).
#
# Most of the edited cells hold plain display text (prices/percentages are
# stored as strings, not numbers), so for any replacement value that Excel
# would otherwise auto-parse as a number we briefly force the cell to Text
# format, write the literal string, then restore the "Normal" style so the
# cell's formatting matches the rest of the untouched column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.376.43'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +0.07%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.841.40'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -0.23%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '0.9992'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  +0.14%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '239.16'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -0.32%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '0.6267'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -0.12%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  +0.11%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.07429'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -0.77%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.2893'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -0.16%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  +1.90%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.07722'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -0.13%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '1.848.42'; ForceText = $false },
    @{ Cell = 'E12'; Value = '  +0.15%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '4.974'; ForceText = $true },
    @{ Cell = 'D14'; Value = '0.6738'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -0.90%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '0.00001027'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -1.98%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '81.76'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -0.34%  '; ForceText = $false },
    @{ Cell = 'E17'; Value = '  +0.41%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '29.434.81'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  +0.18%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '234.10'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  +2.29%  '; ForceText = $false },
    @{ Cell = 'E20'; Value = '  -0.11%  '; ForceText = $false },
    @{ Cell = 'E21'; Value = '  +0.19%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '7.291'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -2.69%  '; ForceText = $false },
    @{ Cell = 'E23'; Value = '  +0.16%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '157.85'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -0.37%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '8.501'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +0.86%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '0.1344'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -1.79%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  -1.27%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '0.07219'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +11.05%  '; ForceText = $false },
    @{ Cell = 'E29'; Value = '  +3.85%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '1.481'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  +0.16%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '4.039'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -1.16%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '4.030'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -1.53%  '; ForceText = $false },
    @{ Cell = 'E33'; Value = '  -0.72%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '1.139'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -0.23%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '0.6973'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  +0.20%  '; ForceText = $false },
    @{ Cell = 'E36'; Value = '  -0.14%  '; ForceText = $false },
    @{ Cell = 'E37'; Value = '  +0.42%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '6.933'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +2.57%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '2.818'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -0.58%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '1.234.41'; ForceText = $false },
    @{ Cell = 'E40'; Value = '  -2.27%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '0.9626'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  +4.74%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '1.000'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  +0.17%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '1.999.51'; ForceText = $false },
    @{ Cell = 'E43'; Value = '  -0.43%  '; ForceText = $false },
    @{ Cell = 'E44'; Value = '  -0.49%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '65.36'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -1.26%  '; ForceText = $false },
    @{ Cell = 'E46'; Value = '  +1.09%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '1.717'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -0.42%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '6.942'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -1.95%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '8.883'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -1.11%  '; ForceText = $false },
    @{ Cell = 'B50'; Value = 'TheSandbox'; ForceText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false },
    @{ Cell = 'D50'; Value = '0.3899'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -1.42%  '; ForceText = $false },
    @{ Cell = 'B51'; Value = 'Algorand'; ForceText = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.1131'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -2.66%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
